$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 <-> Row 39 (Station63): swap species/common/category/reads values
$ws.Range("A38").Value = "Unassigned"
$ws.Range("B38").Value = "Unassigned"
$ws.Range("C38").Value = "Unassigned"
$ws.Range("E38").Value = 298

$ws.Range("A39").Value = "Urophycis sp"
$ws.Range("B39").Value = "Red White or Spotted hake"
$ws.Range("C39").Value = "Teleost Fish"
$ws.Range("E39").Value = 1025

# Row 47 <-> Row 48 (Station69): swap species/common/category/reads values
$ws.Range("A47").Value = "Unassigned"
$ws.Range("B47").Value = "Unassigned"
$ws.Range("C47").Value = "Unassigned"
$ws.Range("E47").Value = 239

$ws.Range("A48").Value = "Urophycis sp"
$ws.Range("B48").Value = "Red White or Spotted hake"
$ws.Range("C48").Value = "Teleost Fish"
$ws.Range("E48").Value = 11
